$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (label unchanged)
$ws.Range("B2").Value = 0.4617684923563146
$ws.Range("C2").Value = 0.4617684923563146
$ws.Range("D2").Value = 0.4617684923563146

# Row 3 - RandomForestRegressor (label unchanged)
$ws.Range("B3").Value = 0.9736619800583912
$ws.Range("C3").Value = 0.9735190221601109
$ws.Range("D3").Value = 0.9723496566305139

# Row 4 - label change: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.979337566335427
$ws.Range("C4").Value = 0.9795312667123779
$ws.Range("D4").Value = 0.978821963618397

# Row 5 - label change: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8467382909027327
$ws.Range("C5").Value = 0.818862000441681
$ws.Range("D5").Value = 0.8210698750958116
